# RotJ - Level 3 up to boss, 560 ahead. Also includes a 30 frames slow boss fight attempt.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FrameCounts")
$ws.Activate()

# --- Row 53 / 54: corrected measurements (B values shift by -1 frame) ---
$ws.Range("B53").Value = 7144
$ws.Range("C53").Value = 7460
$ws.Range("B54").Value = 7174
$ws.Range("C54").Value = 7490

# --- New rows 55-63: frame data through Level 3 boss fight ---
$newRows = @(
    @(55, "X = 282",                  7331, 7647),
    @(56, "Batman appears screen 2",  8507, 8823),
    @(57, "X = 176",                  8588, 8912),
    @(58, "X = 299",                  8648, 8974),
    @(59, "Begin walljump",           9079, 9639),
    @(60, "Black screen",             9164, 9724),
    @(61, "HP = 26",                  9320, 9857),
    @(62, "HP = 0",                   9468, 9997),
    @(63, "Batman disappears",        9810, 10340)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $label = $row[1]
    $b = $row[2]
    $c = $row[3]
    $ws.Range("A" + $r).Value = $label
    $ws.Range("B" + $r).Value = $b
    $ws.Range("C" + $r).Value = $c
}

# Column D holds a shared formula (IF(Bn>0,Cn-Bn,0)); re-assert the formula on
# every touched row so the engine re-binds/recalculates it against the
# newly-populated B/C cells instead of keeping a stale cached result.
for ($r = 53; $r -le 63; $r++) {
    $ws.Range("D" + $r).Formula = "=IF(B" + $r + ">0,C" + $r + "-B" + $r + ",0)"
}

$excel.Calculate()

# --- Scroll / selection bookkeeping to mirror the author's view state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 51
$ws.Range("C64").Select()
